$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1495.85
$ws.Range("I28").Value = 1594.6
$ws.Range("J28").Value = 1199.6
$ws.Range("K28").Value = 1594.6
$ws.Range("L28").Value = 1199.6
$ws.Range("M28").Value = -1109.6
$ws.Range("N28").Value = -2169.6
$ws.Range("H94").Value = 2023.7
$ws.Range("I94").Value = 1673.375
$ws.Range("J94").Value = 3425
$ws.Range("K94").Value = 1673.375
$ws.Range("L94").Value = 3425
$ws.Range("M94").Value = -1222.375
$ws.Range("N94").Value = -4327
$ws.Range("H98").Value = 1655.8206
$ws.Range("I98").Value = 1586.8485
$ws.Range("J98").Value = 2035.1666
$ws.Range("K98").Value = 1586.8485
$ws.Range("L98").Value = 2035.1666
$ws.Range("M98").Value = -88.84850000000006
$ws.Range("N98").Value = -5031.1666
$ws.Range("H112").Value = 20630.936
$ws.Range("J112").Value = 23479.56
$ws.Range("L112").Value = 70438.68000000001
$ws.Range("N112").Value = -72654.68000000001
$ws.Range("H115").Value = 678.8570999999999
$ws.Range("I115").Value = 678.8570999999999
$ws.Range("K115").Value = 2036.5713
$ws.Range("M115").Value = -469.5712999999998
$ws.Range("H121").Value = 3499.25
$ws.Range("J121").Value = 3499.25
$ws.Range("L121").Value = 10497.75
$ws.Range("N121").Value = -13991.75
$ws.Range("H122").Value = 1655.8206
$ws.Range("I122").Value = 1586.8485
$ws.Range("J122").Value = 2035.1666
$ws.Range("K122").Value = 4760.5455
$ws.Range("L122").Value = 6105.4998
$ws.Range("M122").Value = -2310.5455
$ws.Range("N122").Value = -11005.4998
$ws.Range("H132").Value = 1896.2444
$ws.Range("I132").Value = 1661.8948
$ws.Range("J132").Value = 3168.4285
$ws.Range("K132").Value = 4985.6844
$ws.Range("L132").Value = 9505.2855
$ws.Range("M132").Value = -2455.6844
$ws.Range("N132").Value = -14565.2855
$ws.Range("H138").Value = 2414.75
$ws.Range("I138").Value = 785.74194
$ws.Range("J138").Value = 3945.0303
$ws.Range("K138").Value = 2357.22582
$ws.Range("L138").Value = 11835.0909
$ws.Range("M138").Value = 2782.77418
$ws.Range("N138").Value = -22115.0909
$ws.Range("H141").Value = 3786.0833
$ws.Range("I141").Value = 3614.3635
$ws.Range("K141").Value = 10843.0905
$ws.Range("M141").Value = -5663.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12583.207
$ws.Range("I32").Value = 8270.786
$ws.Range("K32").Value = 8270.786
$ws.Range("M32").Value = -7983.786
$ws.Range("H61").Value = 24888.834
$ws.Range("I61").Value = 1908.6666
$ws.Range("K61").Value = 1908.6666
$ws.Range("M61").Value = -1696.6666
$ws.Range("H63").Value = 2521
$ws.Range("I63").Value = 2521
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2521
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1835
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2521
$ws.Range("I66").Value = 2521
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12605
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -9173
$ws.Range("N66").ClearContents()
$ws.Range("H76").Value = 272330.66
$ws.Range("J76").Value = 272330.66
$ws.Range("L76").Value = 272330.66
$ws.Range("N76").Value = -273006.66
$ws.Range("H79").Value = 272330.66
$ws.Range("J79").Value = 272330.66
$ws.Range("L79").Value = 272330.66
$ws.Range("N79").Value = -274670.66
$ws.Range("H80").Value = 30050.334
$ws.Range("J80").Value = 40075.5
$ws.Range("L80").Value = 40075.5
$ws.Range("N80").Value = -42071.5
$ws.Range("H83").Value = 30050.334
$ws.Range("J83").Value = 40075.5
$ws.Range("L83").Value = 120226.5
$ws.Range("N83").Value = -130210.5
$ws.Range("H122").Value = 129903.29
$ws.Range("I122").Value = 1041
$ws.Range("J122").Value = 226550
$ws.Range("K122").Value = 3123
$ws.Range("L122").Value = 679650
$ws.Range("M122").Value = -673
$ws.Range("N122").Value = -684550
$ws.Range("H132").Value = 1658
$ws.Range("I132").Value = 1773
$ws.Range("K132").Value = 5319
$ws.Range("M132").Value = -2789
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -94060
$ws.Range("H136").Value = 24888.834
$ws.Range("I136").Value = 1908.6666
$ws.Range("K136").Value = 5725.9998
$ws.Range("M136").Value = -3175.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 65783.42999999999
$ws.Range("I86").Value = 1622.3334
$ws.Range("J86").Value = 450750
$ws.Range("K86").Value = 1622.3334
$ws.Range("L86").Value = 450750
$ws.Range("M86").Value = -499.3334
$ws.Range("N86").Value = -452996
$ws.Range("H89").Value = 65783.42999999999
$ws.Range("I89").Value = 1622.3334
$ws.Range("J89").Value = 450750
$ws.Range("K89").Value = 8111.666999999999
$ws.Range("L89").Value = 2253750
$ws.Range("M89").Value = -2495.666999999999
$ws.Range("N89").Value = -2264982
$ws.Range("H99").Value = 1608.75
$ws.Range("I99").Value = 1112.5333
$ws.Range("J99").Value = 3097.4
$ws.Range("K99").Value = 1112.5333
$ws.Range("L99").Value = 3097.4
$ws.Range("M99").Value = 385.4666999999999
$ws.Range("N99").Value = -6093.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1862.5588
$ws.Range("I31").Value = 1551.9166
$ws.Range("K31").Value = 1551.9166
$ws.Range("M31").Value = -1256.9166
$ws.Range("H34").Value = 1862.5588
$ws.Range("I34").Value = 1551.9166
$ws.Range("K34").Value = 1551.9166
$ws.Range("M34").Value = -1349.9166
$ws.Range("H69").Value = 37618.5
$ws.Range("I69").Value = 20182
$ws.Range("J69").Value = 55055
$ws.Range("K69").Value = 20182
$ws.Range("L69").Value = 55055
$ws.Range("M69").Value = -19433
$ws.Range("N69").Value = -56553
$ws.Range("H72").Value = 37618.5
$ws.Range("I72").Value = 20182
$ws.Range("J72").Value = 55055
$ws.Range("K72").Value = 60546
$ws.Range("L72").Value = 165165
$ws.Range("M72").Value = -56802
$ws.Range("N72").Value = -172653
$ws.Range("H103").Value = 9500
$ws.Range("I103").Value = 9500
$ws.Range("K103").Value = 9500
$ws.Range("M103").Value = -8328
$ws.Range("H107").Value = 47391.2
$ws.Range("I107").Value = 91137.7
$ws.Range("J107").Value = 3644.7
$ws.Range("K107").Value = 91137.7
$ws.Range("L107").Value = 3644.7
$ws.Range("M107").Value = -89217.7
$ws.Range("N107").Value = -7484.7
$ws.Range("H122").Value = 2005.6923
$ws.Range("I122").Value = 1842
$ws.Range("J122").Value = 2267.6
$ws.Range("K122").Value = 5526
$ws.Range("L122").Value = 6802.799999999999
$ws.Range("M122").Value = -3076
$ws.Range("N122").Value = -11702.8
$ws.Range("H134").Value = 2483.4211
$ws.Range("I134").Value = 1755.0769
$ws.Range("K134").Value = 5265.2307
$ws.Range("M134").Value = -2730.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 377.57144
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = 33
$ws.Range("M15").Value = 107
$ws.Range("H87").Value = 17362.666
$ws.Range("I87").Value = 16120
$ws.Range("K87").Value = 48360
$ws.Range("M87").Value = -47112
$ws.Range("H90").Value = 17362.666
$ws.Range("I90").Value = 16120
$ws.Range("K90").Value = 145080
$ws.Range("M90").Value = -138840
$ws.Range("H132").Value = 2446.353
$ws.Range("I132").Value = 1305.8667
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 11752.8003
$ws.Range("L132").Value = 99000
$ws.Range("M132").Value = -9222.800300000001
$ws.Range("N132").Value = -104060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 20184.193
$ws.Range("I102").Value = 26387.375
$ws.Range("J102").Value = 5588.4707
$ws.Range("K102").Value = 26387.375
$ws.Range("L102").Value = 5588.4707
$ws.Range("M102").Value = -24765.375
$ws.Range("N102").Value = -8832.4707
$ws.Range("H132").Value = 3127.5
$ws.Range("I132").Value = 3118.742
$ws.Range("K132").Value = 9356.226000000001
$ws.Range("M132").Value = -6826.226000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5714.9688
$ws.Range("I22").Value = 2270
$ws.Range("J22").Value = 6207.107
$ws.Range("K22").Value = 2270
$ws.Range("L22").Value = 6207.107
$ws.Range("M22").Value = -1975
$ws.Range("N22").Value = -6797.107
$ws.Range("H27").Value = 5714.9688
$ws.Range("I27").Value = 2270
$ws.Range("J27").Value = 6207.107
$ws.Range("K27").Value = 2270
$ws.Range("L27").Value = 6207.107
$ws.Range("M27").Value = -2163
$ws.Range("N27").Value = -6421.107
$ws.Range("H46").Value = 1635.2609
$ws.Range("I46").Value = 612.1429000000001
$ws.Range("J46").Value = 2082.875
$ws.Range("K46").Value = 612.1429000000001
$ws.Range("L46").Value = 2082.875
$ws.Range("M46").Value = -424.1429000000001
$ws.Range("N46").Value = -2458.875
$ws.Range("H122").Value = 3318.3057
$ws.Range("I122").Value = 3284.0881
$ws.Range("K122").Value = 9852.264299999999
$ws.Range("M122").Value = -7402.264299999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 971.5333000000001
$ws.Range("I107").Value = 610.1905
$ws.Range("K107").Value = 1830.5715
$ws.Range("M107").Value = 89.42849999999999
$ws.Range("H113").Value = 1456.6364
$ws.Range("I113").Value = 1306.6666
$ws.Range("K113").Value = 3919.9998
$ws.Range("M113").Value = -1749.9998
$ws.Range("H122").Value = 2650.8572
$ws.Range("I122").Value = 2797.077
$ws.Range("K122").Value = 8391.231
$ws.Range("M122").Value = -5941.231
$ws.Range("H123").Value = 99997.5
$ws.Range("J123").Value = 99997.5
$ws.Range("L123").Value = 99997.5
$ws.Range("N123").Value = -109797.5
$ws.Range("H125").Value = 57115.832
$ws.Range("J125").Value = 57115.832
$ws.Range("L125").Value = 57115.832
$ws.Range("N125").Value = -66955.83199999999
